# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet right after "总计" and before "2022-Q3"
#    (created by copying the existing "2022-Q3" sheet so it inherits the
#    same header/row formatting), then overwrite its data with the new
#    2022-Q4 fund holdings.
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q4 and shift the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet by duplicating "2022-Q3" so all
# styles/borders/column layout match the existing quarterly sheets.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $totalSheet)

$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Remove the extra rows the copied template had (old "2022-Q3" had 10 data
# rows, rows 8-11; the new "2022-Q4" sheet only needs 6 data rows).
$q4Sheet.Range("A8:H11").Clear()

# Columns B-G hold numeric-looking text (fund codes, percentages, amounts)
# that must stay text, matching the source data's formatting - a leading
# apostrophe forces text storage without leaving a stray number format
# behind.
$q4Sheet.Cells.Item(2, 1).Value = 0
$q4Sheet.Cells.Item(2, 2).Value = "'502000"
$q4Sheet.Cells.Item(2, 3).Value = "'西部利得中证500指数增强（LOF）A"
$q4Sheet.Cells.Item(2, 4).Value = "'18.20"
$q4Sheet.Cells.Item(2, 5).Value = "'90.25"
$q4Sheet.Cells.Item(2, 6).Value = "'2.10"
$q4Sheet.Cells.Item(2, 7).Value = "'0.3822"
$q4Sheet.Cells.Item(2, 8).Value = 3

$q4Sheet.Cells.Item(3, 1).Value = 1
$q4Sheet.Cells.Item(3, 2).Value = "'009300"
$q4Sheet.Cells.Item(3, 3).Value = "'西部利得中证500指数增强（LOF）C"
$q4Sheet.Cells.Item(3, 4).Value = "'4.99"
$q4Sheet.Cells.Item(3, 5).Value = "'90.25"
$q4Sheet.Cells.Item(3, 6).Value = "'2.10"
$q4Sheet.Cells.Item(3, 7).Value = "'0.1048"
$q4Sheet.Cells.Item(3, 8).Value = 3

$q4Sheet.Cells.Item(4, 1).Value = 2
$q4Sheet.Cells.Item(4, 2).Value = "'006441"
$q4Sheet.Cells.Item(4, 3).Value = "'中信建投中证500指数增强C"
$q4Sheet.Cells.Item(4, 4).Value = "'2.27"
$q4Sheet.Cells.Item(4, 5).Value = "'93.50"
$q4Sheet.Cells.Item(4, 6).Value = "'0.83"
$q4Sheet.Cells.Item(4, 7).Value = "'0.0188"
$q4Sheet.Cells.Item(4, 8).Value = 7

$q4Sheet.Cells.Item(5, 1).Value = 3
$q4Sheet.Cells.Item(5, 2).Value = "'006440"
$q4Sheet.Cells.Item(5, 3).Value = "'中信建投中证500指数增强A"
$q4Sheet.Cells.Item(5, 4).Value = "'2.19"
$q4Sheet.Cells.Item(5, 5).Value = "'93.50"
$q4Sheet.Cells.Item(5, 6).Value = "'0.83"
$q4Sheet.Cells.Item(5, 7).Value = "'0.0182"
$q4Sheet.Cells.Item(5, 8).Value = 7

$q4Sheet.Cells.Item(6, 1).Value = 4
$q4Sheet.Cells.Item(6, 2).Value = "'005966"
$q4Sheet.Cells.Item(6, 3).Value = "'安信中证500指数增强C"
$q4Sheet.Cells.Item(6, 4).Value = "'0.16"
$q4Sheet.Cells.Item(6, 5).Value = "'88.79"
$q4Sheet.Cells.Item(6, 6).Value = "'0.89"
$q4Sheet.Cells.Item(6, 7).Value = "'0.0014"
$q4Sheet.Cells.Item(6, 8).Value = 5

$q4Sheet.Cells.Item(7, 1).Value = 5
$q4Sheet.Cells.Item(7, 2).Value = "'005965"
$q4Sheet.Cells.Item(7, 3).Value = "'安信中证500指数增强A"
$q4Sheet.Cells.Item(7, 4).Value = "'0.12"
$q4Sheet.Cells.Item(7, 5).Value = "'88.79"
$q4Sheet.Cells.Item(7, 6).Value = "'0.89"
$q4Sheet.Cells.Item(7, 7).Value = "'0.0011"
$q4Sheet.Cells.Item(7, 8).Value = 5

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary table - insert the 2022-Q4 row on top
# and push the older quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet.Range("A5:D5").Copy()
$totalSheet.Range("A6:D6").PasteSpecial(-4122)

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.53

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(3, 3).Value = 10
$totalSheet.Cells.Item(3, 4).Value = 1.23

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0.4

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(5, 3).Value = 2
$totalSheet.Cells.Item(5, 4).Value = 0

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(6, 3).Value = 3
$totalSheet.Cells.Item(6, 4).Value = 0.24

# ---------------------------------------------------------------------
# Restore the originally active sheet/tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
